$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 534, shifting existing rows 534:606 down to 535:607
$ws.Rows(534).Insert()

# Populate the newly inserted row with the new weekly price observation
$ws.Range("A534").Value = 3
$ws.Range("B534").Value = "Femacal de La Calera"
$ws.Range("C534").Value = "Coquimbo"
$ws.Range("D534").Value = "2023-07-24"
$ws.Range("E534").Value = 5
$ws.Range("F534").Value = 100112043
$ws.Range("G534").Value = "Pepino ensalada"
$ws.Range("H534").Value = "Sin especificar"
$ws.Range("I534").Value = "Primera"
$ws.Range("J534").Value = 100
$ws.Range("K534").Value = 9000
$ws.Range("L534").Value = 10000
$ws.Range("M534").Value = 9500
$ws.Range("N534").Value = "$/caja 60 unidades"
$ws.Range("O534").Value = "Región de Arica y Parinacota"
$ws.Range("P534").Value = 158
$ws.Range("Q534").Value = 60
$ws.Range("R534").Value = "Hortaliza"
